# Apply cryptos list update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $NewValue)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $NewValue
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "43.205.09"
Set-TextValue $ws.Range("E2") "  +0.32%  "

Set-TextValue $ws.Range("D3") "2.323.32"
Set-TextValue $ws.Range("E3") "  +0.93%  "

Set-TextValue $ws.Range("E4") "  +0.00%  "

Set-TextValue $ws.Range("D5") "302.60"
Set-TextValue $ws.Range("E5") "  +0.06%  "

Set-TextValue $ws.Range("D6") "99.67"
Set-TextValue $ws.Range("E6") "  +0.26%  "

Set-TextValue $ws.Range("D7") "0.508"
Set-TextValue $ws.Range("E7") "  +0.14%  "

Set-TextValue $ws.Range("E8") "  +0.04%  "

Set-TextValue $ws.Range("E9") "  +1.73%  "

Set-TextValue $ws.Range("D10") "36.25"
Set-TextValue $ws.Range("E10") "  +5.20%  "

Set-TextValue $ws.Range("E11") "  -0.60%  "

Set-TextValue $ws.Range("E12") "  -0.93%  "

Set-TextValue $ws.Range("D13") "17.63"
Set-TextValue $ws.Range("E13") "  -0.40%  "

Set-TextValue $ws.Range("E14") "  +1.86%  "

Set-TextValue $ws.Range("D15") "2.685.18"
Set-TextValue $ws.Range("E15") "  +0.86%  "

Set-TextValue $ws.Range("D16") "2.344.70"
Set-TextValue $ws.Range("E16") "  +4.19%  "

Set-TextValue $ws.Range("E17") "  -1.10%  "

Set-TextValue $ws.Range("D18") "43.130.89"
Set-TextValue $ws.Range("E18") "  +0.40%  "

Set-TextValue $ws.Range("D19") "12.69"
Set-TextValue $ws.Range("E19") "  +2.80%  "

Set-TextValue $ws.Range("D20") "6.22"
Set-TextValue $ws.Range("E20") "  +1.78%  "

Set-TextValue $ws.Range("E21") "  +0.36%  "

Set-TextValue $ws.Range("D22") "68.25"
Set-TextValue $ws.Range("E22") "  +0.55%  "

Set-TextValue $ws.Range("D23") "241.56"
Set-TextValue $ws.Range("E23") "  +1.87%  "

Set-TextValue $ws.Range("E24") "  -1.51%  "

Set-TextValue $ws.Range("E25") "  -0.42%  "

Set-TextValue $ws.Range("E26") "  -0.05%  "

Set-TextValue $ws.Range("E27") "  +3.80%  "

Set-TextValue $ws.Range("D28") "168.63"
Set-TextValue $ws.Range("E28") "  +0.23%  "

Set-TextValue $ws.Range("E29") "  +1.53%  "

Set-TextValue $ws.Range("E30") "  +0.22%  "

Set-TextValue $ws.Range("D31") "2.03"
Set-TextValue $ws.Range("E31") "  -2.71%  "

Set-TextValue $ws.Range("E32") "  +3.12%  "

Set-TextValue $ws.Range("E33") "  -0.08%  "

Set-TextValue $ws.Range("B34") "RenderToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D34") "4.73"
Set-TextValue $ws.Range("E34") "  +3.72%  "

Set-TextValue $ws.Range("B35") "Celestia"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D35") "17.77"
Set-TextValue $ws.Range("E35") "  +3.93%  "

Set-TextValue $ws.Range("E36") "  -0.99%  "

Set-TextValue $ws.Range("D37") "0.0697"
Set-TextValue $ws.Range("E37") "  -0.35%  "

Set-TextValue $ws.Range("E38") "  +0.68%  "

Set-TextValue $ws.Range("E39") "  +0.52%  "

Set-TextValue $ws.Range("E40") "  -2.11%  "

Set-TextValue $ws.Range("E41") "  +0.26%  "

Set-TextValue $ws.Range("D42") "2.004.69"
Set-TextValue $ws.Range("E42") "  +0.26%  "

Set-TextValue $ws.Range("E43") "  +1.69%  "

Set-TextValue $ws.Range("E44") "  -4.60%  "

Set-TextValue $ws.Range("D45") "10.11"
Set-TextValue $ws.Range("E45") "  -0.16%  "

Set-TextValue $ws.Range("D46") "17.75"
Set-TextValue $ws.Range("E46") "  +0.03%  "

Set-TextValue $ws.Range("E47") "  +0.62%  "

Set-TextValue $ws.Range("B48") "MultiversX"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D48") "55.00"
Set-TextValue $ws.Range("E48") "  -1.35%  "

Set-TextValue $ws.Range("B49") "BitcoinSV"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D49") "76.07"
Set-TextValue $ws.Range("E49") "  +8.32%  "

Set-TextValue $ws.Range("D50") "2.549.74"
Set-TextValue $ws.Range("E50") "  +0.92%  "

Set-TextValue $ws.Range("E51") "  +2.41%  "

